$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the existing data rows (2-12, columns A-D) into memory, keeping the
# source row number as the LAST element (sorting on the array's first
# element is what this host's Sort-Object reliably supports).
$rows = @()
for ($r = 2; $r -le 12; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $d = $ws.Cells.Item($r, 4).Value2
    $rows += ,@($a, $b, $c, $d, $r)
}

# Sort the rows by column A (time) ascending
$sorted = $rows | Sort-Object { $_[0] }

# Write the sorted rows back only where the row actually needs to change
for ($i = 0; $i -lt $sorted.Count; $i++) {
    $destRow = 2 + $i
    $row = $sorted[$i]
    $srcRow = $row[4]
    if ($destRow -ne $srcRow) {
        $ws.Cells.Item($destRow, 1).Value2 = $row[0]
        $ws.Cells.Item($destRow, 2).Value2 = $row[1]
        $ws.Cells.Item($destRow, 3).Value2 = $row[2]
        $ws.Cells.Item($destRow, 4).Value2 = $row[3]
    }
}
